$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B -- shifts old B:K data to C:L,
# matching the price tracker pattern of prepending the newest price-check
# column ahead of all older timestamped columns.
$ws.Range("B1").EntireColumn.Insert()

# Excel assigns the freshly-inserted column the sheet default width; restore
# the tracker column width (internal stored width 21, i.e. ColumnWidth 20.17).
$ws.Range("B1").ColumnWidth = 20.17

# New timestamp header for the freshly scraped column (2025-12-21 12:27 IST)
$ws.Range("B1").Value = "2025-12-21 12:27"

# Latest price snapshot for every SKU row
$ws.Range("B2").Value = 929
$ws.Range("B3").Value = 569
$ws.Range("B4").Value = 299
$ws.Range("B5").Value = 569
$ws.Range("B6").Value = 499
$ws.Range("B7").Value = 569
$ws.Range("B8").Value = 929
$ws.Range("B9").Value = 299
$ws.Range("B10").Value = 299
$ws.Range("B11").Value = 2997
$ws.Range("B12").Value = 569
$ws.Range("B13").Value = 569
$ws.Range("B14").Value = 499
$ws.Range("B15").Value = 499
$ws.Range("B16").Value = 299
$ws.Range("B17").Value = 929
$ws.Range("B18").Value = 499
$ws.Range("B19").Value = 1497
$ws.Range("B20").Value = 929
$ws.Range("B21").Value = 499
$ws.Range("B22").Value = 299
$ws.Range("B23").Value = 1299
$ws.Range("B24").Value = 929
$ws.Range("B25").Value = 929
$ws.Range("B26").Value = 1299
